$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header cell A1
$ws.Range("A1").Value = "VIPIN"

# Update ID values in column A
$ws.Range("A2").Value = 3456
$ws.Range("A3").Value = 5634
$ws.Range("A4").Value = 2345
$ws.Range("A5").Value = 2134

# Add new names in column B, rows 12-16
$ws.Range("B12").Value = "ANU"
$ws.Range("B13").Value = "MEERA"
$ws.Range("B14").Value = "MANU"
$ws.Range("B15").Value = "SNEHA"
$ws.Range("B16").Value = "VIPIN"

# Update selection to match final state
$ws.Range("B17").Select()
